# Update the "Phe khach moi" (New customer fee) pricing table with the
# latest paid-money values (Thu ve, chot, be), replacing the old formula-
# derived figures in B2:R18 with plain literal values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 10
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 13
$ws.Cells.Item(2, 5).Value = 13
$ws.Cells.Item(2, 6).Value = 15
$ws.Cells.Item(2, 7).Value = 15
$ws.Cells.Item(2, 8).Value = 15.5
$ws.Cells.Item(2, 9).Value = 16
$ws.Cells.Item(2, 10).Value = 18
$ws.Cells.Item(2, 11).Value = 19
$ws.Cells.Item(2, 12).Value = 19
$ws.Cells.Item(2, 13).Value = 20
$ws.Cells.Item(2, 14).Value = 20
$ws.Cells.Item(2, 15).Value = 20
$ws.Cells.Item(2, 16).Value = 20
$ws.Cells.Item(2, 17).Value = 20
$ws.Cells.Item(2, 18).Value = 20

# Row 3
$ws.Cells.Item(3, 2).Value = 20
$ws.Cells.Item(3, 3).Value = 24
$ws.Cells.Item(3, 4).Value = 26
$ws.Cells.Item(3, 5).Value = 26
$ws.Cells.Item(3, 6).Value = 30
$ws.Cells.Item(3, 7).Value = 30
$ws.Cells.Item(3, 8).Value = 31
$ws.Cells.Item(3, 9).Value = 32
$ws.Cells.Item(3, 10).Value = 36
$ws.Cells.Item(3, 11).Value = 38
$ws.Cells.Item(3, 12).Value = 38
$ws.Cells.Item(3, 13).Value = 40
$ws.Cells.Item(3, 14).Value = 40
$ws.Cells.Item(3, 15).Value = 40
$ws.Cells.Item(3, 16).Value = 40
$ws.Cells.Item(3, 17).Value = 40
$ws.Cells.Item(3, 18).Value = 40

# Row 4
$ws.Cells.Item(4, 2).Value = 25
$ws.Cells.Item(4, 3).Value = 30
$ws.Cells.Item(4, 4).Value = 32
$ws.Cells.Item(4, 5).Value = 32
$ws.Cells.Item(4, 6).Value = 37
$ws.Cells.Item(4, 7).Value = 37
$ws.Cells.Item(4, 8).Value = 38
$ws.Cells.Item(4, 9).Value = 40
$ws.Cells.Item(4, 10).Value = 44.999999999999993
$ws.Cells.Item(4, 11).Value = 47
$ws.Cells.Item(4, 12).Value = 47
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 50
$ws.Cells.Item(4, 15).Value = 50
$ws.Cells.Item(4, 16).Value = 50
$ws.Cells.Item(4, 17).Value = 50
$ws.Cells.Item(4, 18).Value = 50

# Row 5
$ws.Cells.Item(5, 2).Value = 30
$ws.Cells.Item(5, 3).Value = 36
$ws.Cells.Item(5, 4).Value = 39
$ws.Cells.Item(5, 5).Value = 39
$ws.Cells.Item(5, 6).Value = 45
$ws.Cells.Item(5, 7).Value = 45
$ws.Cells.Item(5, 8).Value = 46
$ws.Cells.Item(5, 9).Value = 48
$ws.Cells.Item(5, 10).Value = 54
$ws.Cells.Item(5, 11).Value = 57
$ws.Cells.Item(5, 12).Value = 57
$ws.Cells.Item(5, 13).Value = 60
$ws.Cells.Item(5, 14).Value = 60
$ws.Cells.Item(5, 15).Value = 60
$ws.Cells.Item(5, 16).Value = 60
$ws.Cells.Item(5, 17).Value = 60
$ws.Cells.Item(5, 18).Value = 60

# Row 6
$ws.Cells.Item(6, 2).Value = 35
$ws.Cells.Item(6, 3).Value = 42
$ws.Cells.Item(6, 4).Value = 46
$ws.Cells.Item(6, 5).Value = 46
$ws.Cells.Item(6, 6).Value = 53
$ws.Cells.Item(6, 7).Value = 53
$ws.Cells.Item(6, 8).Value = 54
$ws.Cells.Item(6, 9).Value = 56
$ws.Cells.Item(6, 10).Value = 63
$ws.Cells.Item(6, 11).Value = 67
$ws.Cells.Item(6, 12).Value = 67
$ws.Cells.Item(6, 13).Value = 70
$ws.Cells.Item(6, 14).Value = 70
$ws.Cells.Item(6, 15).Value = 70
$ws.Cells.Item(6, 16).Value = 70
$ws.Cells.Item(6, 17).Value = 70
$ws.Cells.Item(6, 18).Value = 70

# Row 7
$ws.Cells.Item(7, 2).Value = 40
$ws.Cells.Item(7, 3).Value = 48
$ws.Cells.Item(7, 4).Value = 52
$ws.Cells.Item(7, 5).Value = 52
$ws.Cells.Item(7, 6).Value = 60
$ws.Cells.Item(7, 7).Value = 60
$ws.Cells.Item(7, 8).Value = 62
$ws.Cells.Item(7, 9).Value = 64
$ws.Cells.Item(7, 10).Value = 72
$ws.Cells.Item(7, 11).Value = 76
$ws.Cells.Item(7, 12).Value = 76
$ws.Cells.Item(7, 13).Value = 80
$ws.Cells.Item(7, 14).Value = 80
$ws.Cells.Item(7, 15).Value = 80
$ws.Cells.Item(7, 16).Value = 80
$ws.Cells.Item(7, 17).Value = 80
$ws.Cells.Item(7, 18).Value = 80

# Row 8
$ws.Cells.Item(8, 2).Value = 50
$ws.Cells.Item(8, 3).Value = 60
$ws.Cells.Item(8, 4).Value = 65
$ws.Cells.Item(8, 5).Value = 65
$ws.Cells.Item(8, 6).Value = 75
$ws.Cells.Item(8, 7).Value = 75
$ws.Cells.Item(8, 8).Value = 77
$ws.Cells.Item(8, 9).Value = 80
$ws.Cells.Item(8, 10).Value = 89.999999999999986
$ws.Cells.Item(8, 11).Value = 95
$ws.Cells.Item(8, 12).Value = 95
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 14).Value = 100
$ws.Cells.Item(8, 15).Value = 100
$ws.Cells.Item(8, 16).Value = 100
$ws.Cells.Item(8, 17).Value = 100
$ws.Cells.Item(8, 18).Value = 100

# Row 9
$ws.Cells.Item(9, 2).Value = 60
$ws.Cells.Item(9, 3).Value = 72
$ws.Cells.Item(9, 4).Value = 78
$ws.Cells.Item(9, 5).Value = 78
$ws.Cells.Item(9, 6).Value = 90
$ws.Cells.Item(9, 7).Value = 90
$ws.Cells.Item(9, 8).Value = 93
$ws.Cells.Item(9, 9).Value = 96
$ws.Cells.Item(9, 10).Value = 108
$ws.Cells.Item(9, 11).Value = 114
$ws.Cells.Item(9, 12).Value = 114
$ws.Cells.Item(9, 13).Value = 120
$ws.Cells.Item(9, 14).Value = 120
$ws.Cells.Item(9, 15).Value = 120
$ws.Cells.Item(9, 16).Value = 120
$ws.Cells.Item(9, 17).Value = 120
$ws.Cells.Item(9, 18).Value = 120

# Row 10
$ws.Cells.Item(10, 2).Value = 75
$ws.Cells.Item(10, 3).Value = 90
$ws.Cells.Item(10, 4).Value = 97
$ws.Cells.Item(10, 5).Value = 97
$ws.Cells.Item(10, 6).Value = 112
$ws.Cells.Item(10, 7).Value = 112
$ws.Cells.Item(10, 8).Value = 116
$ws.Cells.Item(10, 9).Value = 120
$ws.Cells.Item(10, 10).Value = 134.99999999999997
$ws.Cells.Item(10, 11).Value = 142
$ws.Cells.Item(10, 12).Value = 142
$ws.Cells.Item(10, 13).Value = 150
$ws.Cells.Item(10, 14).Value = 150
$ws.Cells.Item(10, 15).Value = 150
$ws.Cells.Item(10, 16).Value = 150
$ws.Cells.Item(10, 17).Value = 150
$ws.Cells.Item(10, 18).Value = 150

# Row 11
$ws.Cells.Item(11, 2).Value = 80
$ws.Cells.Item(11, 3).Value = 96
$ws.Cells.Item(11, 4).Value = 104
$ws.Cells.Item(11, 5).Value = 104
$ws.Cells.Item(11, 6).Value = 120
$ws.Cells.Item(11, 7).Value = 120
$ws.Cells.Item(11, 8).Value = 124
$ws.Cells.Item(11, 9).Value = 128
$ws.Cells.Item(11, 10).Value = 144
$ws.Cells.Item(11, 11).Value = 152
$ws.Cells.Item(11, 12).Value = 152
$ws.Cells.Item(11, 13).Value = 160
$ws.Cells.Item(11, 14).Value = 160
$ws.Cells.Item(11, 15).Value = 160
$ws.Cells.Item(11, 16).Value = 160
$ws.Cells.Item(11, 17).Value = 160
$ws.Cells.Item(11, 18).Value = 160

# Row 12
$ws.Cells.Item(12, 2).Value = 100
$ws.Cells.Item(12, 3).Value = 120
$ws.Cells.Item(12, 4).Value = 130
$ws.Cells.Item(12, 5).Value = 130
$ws.Cells.Item(12, 6).Value = 150
$ws.Cells.Item(12, 7).Value = 150
$ws.Cells.Item(12, 8).Value = 155
$ws.Cells.Item(12, 9).Value = 160
$ws.Cells.Item(12, 10).Value = 179.99999999999997
$ws.Cells.Item(12, 11).Value = 190
$ws.Cells.Item(12, 12).Value = 190
$ws.Cells.Item(12, 13).Value = 200
$ws.Cells.Item(12, 14).Value = 200
$ws.Cells.Item(12, 15).Value = 200
$ws.Cells.Item(12, 16).Value = 200
$ws.Cells.Item(12, 17).Value = 200
$ws.Cells.Item(12, 18).Value = 200

# Row 13
$ws.Cells.Item(13, 2).Value = 125
$ws.Cells.Item(13, 3).Value = 150
$ws.Cells.Item(13, 4).Value = 162
$ws.Cells.Item(13, 5).Value = 162
$ws.Cells.Item(13, 6).Value = 187
$ws.Cells.Item(13, 7).Value = 187
$ws.Cells.Item(13, 8).Value = 194
$ws.Cells.Item(13, 9).Value = 200
$ws.Cells.Item(13, 10).Value = 225
$ws.Cells.Item(13, 11).Value = 237
$ws.Cells.Item(13, 12).Value = 237
$ws.Cells.Item(13, 13).Value = 250
$ws.Cells.Item(13, 14).Value = 250
$ws.Cells.Item(13, 15).Value = 250
$ws.Cells.Item(13, 16).Value = 250
$ws.Cells.Item(13, 17).Value = 250
$ws.Cells.Item(13, 18).Value = 250

# Row 14
$ws.Cells.Item(14, 2).Value = 150
$ws.Cells.Item(14, 3).Value = 180
$ws.Cells.Item(14, 4).Value = 195
$ws.Cells.Item(14, 5).Value = 195
$ws.Cells.Item(14, 6).Value = 225
$ws.Cells.Item(14, 7).Value = 225
$ws.Cells.Item(14, 8).Value = 233
$ws.Cells.Item(14, 9).Value = 240
$ws.Cells.Item(14, 10).Value = 269.99999999999994
$ws.Cells.Item(14, 11).Value = 285
$ws.Cells.Item(14, 12).Value = 285
$ws.Cells.Item(14, 13).Value = 300
$ws.Cells.Item(14, 14).Value = 300
$ws.Cells.Item(14, 15).Value = 300
$ws.Cells.Item(14, 16).Value = 300
$ws.Cells.Item(14, 17).Value = 300
$ws.Cells.Item(14, 18).Value = 300

# Row 15
$ws.Cells.Item(15, 2).Value = 175
$ws.Cells.Item(15, 3).Value = 210
$ws.Cells.Item(15, 4).Value = 227
$ws.Cells.Item(15, 5).Value = 227
$ws.Cells.Item(15, 6).Value = 262
$ws.Cells.Item(15, 7).Value = 262
$ws.Cells.Item(15, 8).Value = 271
$ws.Cells.Item(15, 9).Value = 280
$ws.Cells.Item(15, 10).Value = 315
$ws.Cells.Item(15, 11).Value = 380
$ws.Cells.Item(15, 12).Value = 380
$ws.Cells.Item(15, 13).Value = 380
$ws.Cells.Item(15, 14).Value = 400
$ws.Cells.Item(15, 15).Value = 400
$ws.Cells.Item(15, 16).Value = 400
$ws.Cells.Item(15, 17).Value = 400
$ws.Cells.Item(15, 18).Value = 400

# Row 16
$ws.Cells.Item(16, 2).Value = 250
$ws.Cells.Item(16, 3).Value = 300
$ws.Cells.Item(16, 4).Value = 325
$ws.Cells.Item(16, 5).Value = 325
$ws.Cells.Item(16, 6).Value = 375
$ws.Cells.Item(16, 7).Value = 375
$ws.Cells.Item(16, 8).Value = 387.5
$ws.Cells.Item(16, 9).Value = 400
$ws.Cells.Item(16, 10).Value = 450
$ws.Cells.Item(16, 11).Value = 475
$ws.Cells.Item(16, 12).Value = 475
$ws.Cells.Item(16, 13).Value = 500
$ws.Cells.Item(16, 14).Value = 500
$ws.Cells.Item(16, 15).Value = 500
$ws.Cells.Item(16, 16).Value = 500
$ws.Cells.Item(16, 17).Value = 500
$ws.Cells.Item(16, 18).Value = 500

# Row 17
$ws.Cells.Item(17, 2).Value = 300
$ws.Cells.Item(17, 3).Value = 360
$ws.Cells.Item(17, 4).Value = 390
$ws.Cells.Item(17, 5).Value = 390
$ws.Cells.Item(17, 6).Value = 450
$ws.Cells.Item(17, 7).Value = 450
$ws.Cells.Item(17, 8).Value = 465
$ws.Cells.Item(17, 9).Value = 480
$ws.Cells.Item(17, 10).Value = 540
$ws.Cells.Item(17, 11).Value = 570
$ws.Cells.Item(17, 12).Value = 570
$ws.Cells.Item(17, 13).Value = 600
$ws.Cells.Item(17, 14).Value = 600
$ws.Cells.Item(17, 15).Value = 600
$ws.Cells.Item(17, 16).Value = 600
$ws.Cells.Item(17, 17).Value = 600
$ws.Cells.Item(17, 18).Value = 600

# Row 18
$ws.Cells.Item(18, 2).Value = 350
$ws.Cells.Item(18, 3).Value = 420
$ws.Cells.Item(18, 4).Value = 454
$ws.Cells.Item(18, 5).Value = 454
$ws.Cells.Item(18, 6).Value = 524
$ws.Cells.Item(18, 7).Value = 524
$ws.Cells.Item(18, 8).Value = 542
$ws.Cells.Item(18, 9).Value = 560
$ws.Cells.Item(18, 10).Value = 630
$ws.Cells.Item(18, 11).Value = 664
$ws.Cells.Item(18, 12).Value = 664
$ws.Cells.Item(18, 13).Value = 700
$ws.Cells.Item(18, 14).Value = 700
$ws.Cells.Item(18, 15).Value = 700
$ws.Cells.Item(18, 16).Value = 700
$ws.Cells.Item(18, 17).Value = 700
$ws.Cells.Item(18, 18).Value = 700

# Move / persist the active selection to L14, matching the saved sheet view
$ws.Range("L14").Select() | Out-Null
